$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.544073932366012
$ws.Range("C2").Value = 0.482679459634459
$ws.Range("D2").Value = 0.0923565941666421
$ws.Range("E2").Value = 0.899049356917915
$ws.Range("F2").Value = 0.276657739058775
$ws.Range("G2").Value = 0.0958001000676929
$ws.Range("H2").Value = 0.00332577920357889
$ws.Range("I2").Value = 0.994290255172617
$ws.Range("J2").Value = 0.0942990846749272
$ws.Range("K2").Value = 0.871942784825029
$ws.Range("L2").Value = 0.000323748418047503
$ws.Range("M2").Value = 0.018895134944227
$ws.Range("N2").Value = 0.00103010860287842
$ws.Range("O2").Value = 0.000147158371839774
$ws.Range("P2").Value = 0.0000588633487359096
$ws.Range("Q2").Value = 0.998175236189187
$ws.Range("R2").Value = 0.00108897195161433
$ws.Range("S2").Value = 0.973570356417577
$ws.Range("T2").Value = 0.000117726697471819
$ws.Range("U2").Value = 0.000294316743679548
$ws.Range("V2").Value = 0.048856579450805
$ws.Range("W2").Value = 0.0331694970126851
$ws.Range("X2").Value = 0.0461488654089531
$ws.Range("B3").Value = 0.0512405450746093
$ws.Range("C3").Value = 0.280012949936722
$ws.Range("D3").Value = 0.390205138770345
$ws.Range("E3").Value = 0.0251640815846013
$ws.Range("F3").Value = 0.621508667628101
$ws.Range("G3").Value = 0.45934014186067
$ws.Range("H3").Value = 0.000559201812991141
$ws.Range("I3").Value = 0.00188362715954911
$ws.Range("J3").Value = 0.00441475115519322
$ws.Range("K3").Value = 0.000323748418047503
$ws.Range("L3").Value = 0.499955852488448
$ws.Range("M3").Value = 0.9116461135474
$ws.Range("N3").Value = 0.00126556199782206
$ws.Range("O3").Value = 0.000264885069311593
$ws.Range("P3").Value = 0.922771286458487
$ws.Range("Q3").Value = 0.000412043441151367
$ws.Range("R3").Value = 0.0000588633487359096
$ws.Range("S3").Value = 0.0000294316743679548
$ws.Range("T3").Value = 0.999058186420225
$ws.Range("U3").Value = 0.998057509491715
$ws.Range("V3").Value = 0.0116549430497101
$ws.Range("W3").Value = 0.024722606469082
$ws.Range("X3").Value = 0.0978603172734497
$ws.Range("B4").Value = 0.393118874532772
$ws.Range("C4").Value = 0.214115431026871
$ws.Range("D4").Value = 0.363804926862289
$ws.Range("E4").Value = 0.0524178120493275
$ws.Range("F4").Value = 0.0495335079612679
$ws.Range("G4").Value = 0.0522117903287518
$ws.Range("H4").Value = 0.99511434205492
$ws.Range("I4").Value = 0.0012361303234541
$ws.Range("J4").Value = 0.899755717102746
$ws.Range("K4").Value = 0.125879271271743
$ws.Range("L4").Value = 0.00185419548518115
$ws.Range("M4").Value = 0.00258998734438002
$ws.Range("N4").Value = 0.000147158371839774
$ws.Range("O4").Value = 0.000294316743679548
$ws.Range("P4").Value = 0.000235453394943638
$ws.Range("Q4").Value = 0.000971245254142508
$ws.Range("R4").Value = 0.998528416281602
$ws.Range("S4").Value = 0.0262824852105836
$ws.Range("T4").Value = 0.000470906789887277
$ws.Range("U4").Value = 0.000264885069311593
$ws.Range("V4").Value = 0.0340230155693557
$ws.Range("W4").Value = 0.93595667657533
$ws.Range("X4").Value = 0.816876122082585
$ws.Range("B5").Value = 0.0113900579803985
$ws.Range("C5").Value = 0.0229567060070047
$ws.Range("D5").Value = 0.153486181828884
$ws.Range("E5").Value = 0.0233098860994202
$ws.Range("F5").Value = 0.0520352002825441
$ws.Range("G5").Value = 0.392324219324837
$ws.Range("H5").Value = 0.00100067692851046
$ws.Range("I5").Value = 0.00258998734438002
$ws.Range("J5").Value = 0.00141272036966183
$ws.Range("K5").Value = 0.00179533213644524
$ws.Range("L5").Value = 0.497660181887748
$ws.Range("M5").Value = 0.0668099008152574
$ws.Range("N5").Value = 0.997527739353092
$ws.Range("O5").Value = 0.999293639815169
$ws.Range("P5").Value = 0.0769049651234659
$ws.Range("Q5").Value = 0.000441475115519322
$ws.Range("R5").Value = 0.000323748418047503
$ws.Range("S5").Value = 0.0000588633487359096
$ws.Range("T5").Value = 0.000353180092415458
$ws.Range("U5").Value = 0.00135385702092592
$ws.Range("V5").Value = 0.905406598581393
$ws.Range("W5").Value = 0.00615121994290255
$ws.Range("X5").Value = 0.0390264002119081
